$wb = $excel.ActiveWorkbook

# Update "想去人数" (F) and "最低票价" (G) values on the 展览 sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 2107
$ws1.Range("G7").Value = 60
$ws1.Range("F10").Value = 4649
$ws1.Range("F13").Value = 292
$ws1.Range("F15").Value = 19
$ws1.Range("F16").Value = 149
$ws1.Range("F19").Value = 94
$ws1.Range("F20").Value = 3535
$ws1.Range("F25").Value = 91
$ws1.Range("F26").Value = 104
$ws1.Range("F32").Value = 741
$ws1.Range("F33").Value = 2176
$ws1.Range("F34").Value = 402

# Update the same values on the 全部类型 sheet (rows shifted by 1 after row 26)
$ws2 = $wb.Worksheets.Item("全部类型")
$ws2.Range("F7").Value = 2107
$ws2.Range("G7").Value = 60
$ws2.Range("F10").Value = 4649
$ws2.Range("F13").Value = 292
$ws2.Range("F15").Value = 19
$ws2.Range("F16").Value = 149
$ws2.Range("F19").Value = 94
$ws2.Range("F20").Value = 3535
$ws2.Range("F25").Value = 91
$ws2.Range("F26").Value = 104
$ws2.Range("F33").Value = 741
$ws2.Range("F34").Value = 2176
$ws2.Range("F35").Value = 402
